$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Move the hidden "_GoBack" bookmark from around ${name} (near the
#    top of the document) to around ${ort}, ${certda} (near the end,
#    in the signature/date line). Word auto-renumbers every other
#    bookmark's w:id sequentially on save, so simply deleting the old
#    one and adding the new one reproduces the id shift seen in the
#    diff automatically.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$goBackTarget = $d.Content
$null = $goBackTarget.Find.Execute('${ort}, ${certda}')
$d.Bookmarks.Add("_GoBack", $goBackTarget)

# ------------------------------------------------------------------
# 2) The "${informatik}" placeholder used to be split across three
#    runs ("${", "informatik", "}"). Collapse it back into a single
#    run/text node by replacing the text in place.
# ------------------------------------------------------------------
$informatikRange = $d.Content
$informatikRange.Find.Execute('${informatik}', $true, $false, $false, $false, $false, `
    $true, 1, $false, '${informatik}', 2)

Write-Output "done"
